$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fecha (D) and Volumen (J) values for rows 2-4.
# The data was reordered: row2 <- old row4's date logic etc.
$ws.Range("D2").Value = 44692
$ws.Range("J2").Value = 120

$ws.Range("D3").Value = 44691
$ws.Range("J3").Value = 100

$ws.Range("D4").Value = 44687
$ws.Range("J4").Value = 160
